# Apply the "stuff at the bottom of the sheets" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pair_kind ("generic") values for the practice rows (J2:J5)
$ws.Range("J2:J5").Value = "generic"

# New "stim details" block starting at row 27
$ws.Range("A27").Value = "stim details"

$headers = @("month", "word_type", "need_audio", "need_image", "word", "count", "find images")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(28, $i + 1).Value = $headers[$i]
}

$rows = @(
    @(6, "video"),
    @(6, "video"),
    @(7, "video"),
    @(7, "video"),
    @(6, "audio"),
    @(6, "audio"),
    @(7, "audio"),
    @(7, "audio")
)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 29 + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
